# Regenerate the handoff report with a new package GUID / commit hash and
# updated timestamps ("Generate Report for Handoff").

$wb = $excel.ActiveWorkbook

$oldGuid = "ffc070f9-4506-4e44-9c6b-29544f01d669"
$newGuid = "e5102ddf-f96b-4857-b895-2f8760054d0f"

$oldCommit = "941893b560327afbbdc924044f7d72859d084926"
$newCommit = "9dbf474747de209f42181bf3636799be06e68780"

$newMdName = "$newGuid.md"
$newZhName = "$newGuid.$newCommit.zh-cn.xlf"
$newDeName = "$newGuid.$newCommit.de-de.xlf"

function Set-HyperlinkDisplay($ws, $cellAddress, $newText) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $cellAddress) {
            $hl.TextToDisplay = $newText
        }
    }
}

# --- Sheet 1: "Overview" ---
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A2").Value = $newMdName
Set-HyperlinkDisplay $ws1 "`$A`$2" $newMdName

$ws1.Range("D2").Value = "2016-51-18 16:51:17"

# --- Sheet 2: "zh-cn" ---
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A2").Value = $newMdName
Set-HyperlinkDisplay $ws2 "`$A`$2" $newMdName

$ws2.Range("D2").Value = $newZhName
Set-HyperlinkDisplay $ws2 "`$D`$2" $newZhName

$ws2.Range("E2").Value = "2016-03-18 16:51:14"

# --- Sheet 3: "de-de" ---
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A2").Value = $newMdName
Set-HyperlinkDisplay $ws3 "`$A`$2" $newMdName

$ws3.Range("D2").Value = $newDeName
Set-HyperlinkDisplay $ws3 "`$D`$2" $newDeName

$ws3.Range("E2").Value = "2016-03-18 16:51:17"
